# Update the workbook to add a new "May 2020" sheet with the latest
# Honorable Rock monthly record, duplicating the structure/format of the
# most recent existing sheet ("April 2020") and updating its data.

$wb = $excel.ActiveWorkbook

# Locate the most recent existing monthly sheet - this one currently
# holds the formatting / formulas we want to reuse for the new month.
$src = $wb.Worksheets.Item("April 2020")

# Duplicate it (keeps all formulas, styles, column widths, merged cells,
# page setup, etc.) and place the copy right after the source sheet.
$src.Copy([System.Reflection.Missing]::Value, $src)
$ws = $wb.Worksheets.Item($src.Index + 1)
$ws.Name = "May 2020"

# Update the month label used by the JSON-building formula in E2.
$ws.Range("E2").Formula = '="    """&"May 2020"&""""&":"'

# Replace the rank/guild-name/contribution data for the new month. The
# JSON-row formulas in column E (and the rank numbers in column B) are
# unchanged from the copied sheet, so only guild name (C) and
# contribution (D) need new values per ranked row (rank numbers in B
# also rewritten for completeness/safety).
$ws.Cells.Item(4, 2).Value = 1
$ws.Cells.Item(4, 3).Value = 'Savages'
$ws.Cells.Item(4, 4).Value = 269596089
$ws.Cells.Item(5, 2).Value = 2
$ws.Cells.Item(5, 3).Value = 'Smile'
$ws.Cells.Item(5, 4).Value = 267669855
$ws.Cells.Item(6, 2).Value = 3
$ws.Cells.Item(6, 3).Value = 'Eternal'
$ws.Cells.Item(6, 4).Value = 259710759
$ws.Cells.Item(7, 2).Value = 4
$ws.Cells.Item(7, 3).Value = 'Bounce'
$ws.Cells.Item(7, 4).Value = 239326724
$ws.Cells.Item(8, 2).Value = 5
$ws.Cells.Item(8, 3).Value = 'Elite'
$ws.Cells.Item(8, 4).Value = 208546754
$ws.Cells.Item(9, 2).Value = 6
$ws.Cells.Item(9, 3).Value = 'Spring'
$ws.Cells.Item(9, 4).Value = 151062021
$ws.Cells.Item(10, 2).Value = 7
$ws.Cells.Item(10, 3).Value = 'Epic'
$ws.Cells.Item(10, 4).Value = 150511738
$ws.Cells.Item(11, 2).Value = 8
$ws.Cells.Item(11, 3).Value = 'Sunset'
$ws.Cells.Item(11, 4).Value = 147134279
$ws.Cells.Item(12, 2).Value = 9
$ws.Cells.Item(12, 3).Value = 'Beaters'
$ws.Cells.Item(12, 4).Value = 137104314
$ws.Cells.Item(13, 2).Value = 10
$ws.Cells.Item(13, 3).Value = 'Downtime'
$ws.Cells.Item(13, 4).Value = 125989790
$ws.Cells.Item(14, 2).Value = 11
$ws.Cells.Item(14, 3).Value = 'lolicafe'
$ws.Cells.Item(14, 4).Value = 111242227
$ws.Cells.Item(15, 2).Value = 12
$ws.Cells.Item(15, 3).Value = 'Imperium'
$ws.Cells.Item(15, 4).Value = 106048719
$ws.Cells.Item(16, 2).Value = 13
$ws.Cells.Item(16, 3).Value = 'Remorse'
$ws.Cells.Item(16, 4).Value = 99620457
$ws.Cells.Item(17, 2).Value = 14
$ws.Cells.Item(17, 3).Value = 'Revive'
$ws.Cells.Item(17, 4).Value = 96167171
$ws.Cells.Item(18, 2).Value = 15
$ws.Cells.Item(18, 3).Value = 'Cleanse'
$ws.Cells.Item(18, 4).Value = 94661918
$ws.Cells.Item(19, 2).Value = 16
$ws.Cells.Item(19, 3).Value = 'Maha'
$ws.Cells.Item(19, 4).Value = 93179716
$ws.Cells.Item(20, 2).Value = 17
$ws.Cells.Item(20, 3).Value = 'Gintama'
$ws.Cells.Item(20, 4).Value = 92515414
$ws.Cells.Item(21, 2).Value = 18
$ws.Cells.Item(21, 3).Value = 'Lithe'
$ws.Cells.Item(21, 4).Value = 90930849
$ws.Cells.Item(22, 2).Value = 19
$ws.Cells.Item(22, 3).Value = 'Undertale'
$ws.Cells.Item(22, 4).Value = 89532722
$ws.Cells.Item(23, 2).Value = 20
$ws.Cells.Item(23, 3).Value = 'RainSong'
$ws.Cells.Item(23, 4).Value = 87759532
$ws.Cells.Item(24, 2).Value = 21
$ws.Cells.Item(24, 3).Value = 'Broke'
$ws.Cells.Item(24, 4).Value = 86226610
$ws.Cells.Item(25, 2).Value = 22
$ws.Cells.Item(25, 3).Value = 'Erda'
$ws.Cells.Item(25, 4).Value = 85486323
$ws.Cells.Item(26, 2).Value = 23
$ws.Cells.Item(26, 3).Value = 'Sora'
$ws.Cells.Item(26, 4).Value = 80366976
$ws.Cells.Item(27, 2).Value = 24
$ws.Cells.Item(27, 3).Value = 'Sugar'
$ws.Cells.Item(27, 4).Value = 71508173
$ws.Cells.Item(28, 2).Value = 25
$ws.Cells.Item(28, 3).Value = 'Rising'
$ws.Cells.Item(28, 4).Value = 69473960
$ws.Cells.Item(29, 2).Value = 26
$ws.Cells.Item(29, 3).Value = 'Earnest'
$ws.Cells.Item(29, 4).Value = 64630310
$ws.Cells.Item(30, 2).Value = 27
$ws.Cells.Item(30, 3).Value = 'Aloe'
$ws.Cells.Item(30, 4).Value = 63263190
$ws.Cells.Item(31, 2).Value = 28
$ws.Cells.Item(31, 3).Value = 'Howl'
$ws.Cells.Item(31, 4).Value = 63190342
$ws.Cells.Item(32, 2).Value = 29
$ws.Cells.Item(32, 3).Value = 'Oceania'
$ws.Cells.Item(32, 4).Value = 62327891
$ws.Cells.Item(33, 2).Value = 30
$ws.Cells.Item(33, 3).Value = 'Ravers'
$ws.Cells.Item(33, 4).Value = 61025622
$ws.Cells.Item(34, 2).Value = 31
$ws.Cells.Item(34, 3).Value = 'Fabled'
$ws.Cells.Item(34, 4).Value = 57941186
$ws.Cells.Item(35, 2).Value = 32
$ws.Cells.Item(35, 3).Value = 'Mystical'
$ws.Cells.Item(35, 4).Value = 57268486
$ws.Cells.Item(36, 2).Value = 33
$ws.Cells.Item(36, 3).Value = 'RainDrop'
$ws.Cells.Item(36, 4).Value = 55118248
$ws.Cells.Item(37, 2).Value = 34
$ws.Cells.Item(37, 3).Value = 'Atelier'
$ws.Cells.Item(37, 4).Value = 51024516
$ws.Cells.Item(38, 2).Value = 35
$ws.Cells.Item(38, 3).Value = 'Exorcist'
$ws.Cells.Item(38, 4).Value = 49167200
$ws.Cells.Item(39, 2).Value = 36
$ws.Cells.Item(39, 3).Value = 'Nutsy'
$ws.Cells.Item(39, 4).Value = 48971834
$ws.Cells.Item(40, 2).Value = 37
$ws.Cells.Item(40, 3).Value = 'Kingdom'
$ws.Cells.Item(40, 4).Value = 48553741
$ws.Cells.Item(41, 2).Value = 38
$ws.Cells.Item(41, 3).Value = 'Path'
$ws.Cells.Item(41, 4).Value = 46873503
$ws.Cells.Item(42, 2).Value = 39
$ws.Cells.Item(42, 3).Value = 'Tama'
$ws.Cells.Item(42, 4).Value = 46409016
$ws.Cells.Item(43, 2).Value = 40
$ws.Cells.Item(43, 3).Value = 'CyberThreat'
$ws.Cells.Item(43, 4).Value = 45954662
$ws.Cells.Item(44, 2).Value = 41
$ws.Cells.Item(44, 3).Value = 'Fandom'
$ws.Cells.Item(44, 4).Value = 44791209
$ws.Cells.Item(45, 2).Value = 42
$ws.Cells.Item(45, 3).Value = 'Miao'
$ws.Cells.Item(45, 4).Value = 44427832
$ws.Cells.Item(46, 2).Value = 43
$ws.Cells.Item(46, 3).Value = 'Weibo'
$ws.Cells.Item(46, 4).Value = 43921823
$ws.Cells.Item(47, 2).Value = 44
$ws.Cells.Item(47, 3).Value = 'Bubbles'
$ws.Cells.Item(47, 4).Value = 41345760
$ws.Cells.Item(48, 2).Value = 45
$ws.Cells.Item(48, 3).Value = 'Reboot'
$ws.Cells.Item(48, 4).Value = 41174173
$ws.Cells.Item(49, 2).Value = 46
$ws.Cells.Item(49, 3).Value = 'Comity'
$ws.Cells.Item(49, 4).Value = 40557398
$ws.Cells.Item(50, 2).Value = 47
$ws.Cells.Item(50, 3).Value = 'Supahot'
$ws.Cells.Item(50, 4).Value = 40498998
$ws.Cells.Item(51, 2).Value = 48
$ws.Cells.Item(51, 3).Value = 'Faction'
$ws.Cells.Item(51, 4).Value = 39836403
$ws.Cells.Item(52, 2).Value = 49
$ws.Cells.Item(52, 3).Value = 'Artifacts'
$ws.Cells.Item(52, 4).Value = 39468858
$ws.Cells.Item(53, 2).Value = 50
$ws.Cells.Item(53, 3).Value = 'chigga'
$ws.Cells.Item(53, 4).Value = 39158297

# Make the newly added sheet the active/selected tab, matching the
# book's active-tab bookkeeping (this also clears the "selected" flag
# on the previously active "April 2020" sheet), and put the selection
# on the freshly edited header formula cell (E2) rather than wherever
# the copied sheet's selection happened to be.
$ws.Select()
$ws.Range("E2").Select()
